# disability_prevalence.xlsx (Tianeti) update:
# - Replace the placeholder "confidential" block (rows 4-5) with two real
#   data rows: "family with disabilities Persons" and "disabilities Persons".
# - Title (row 1) text changes and now spans A1:I1.
# - Row 6 becomes the merged "Source:" attribution line (A6:H6).
# - Column A gets a Sylfaen font as its base/default font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column A: base/default font for the column ----
$ws.Columns.Item(1).Font.Name = "Sylfaen"
$ws.Columns.Item(1).Font.Size = 11
$ws.Columns.Item(1).ColumnWidth = 19.95

# ---- Row 1: title (merged A1:I1) ----
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Tianeti Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 51

# ---- Row 2: subtitle, unchanged text/style - just restore default height ----
$ws.Rows.Item(2).AutoFit()

# ---- Row 3: blank cell under title keeps its top border but now uses the
#      column's Sylfaen font (no explicit font override) ----
$ws.Range("A3").ClearContents()

# ---- Row 4: "family with disabilities Persons" + real data ----
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.Color = -16777216
$ws.Range("A4").Font.Bold = $false
$ws.Range("A4").Font.Underline = 0
$ws.Range("A4").Interior.ThemeColor = 0
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(8).Weight = 2
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 24.75

$data4 = @(244, 215, 208, 215, 217, 221, 229, 229)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt 8; $i++) {
    $addr = "$($cols[$i])4"
    $ws.Range($addr).Value = $data4[$i]
    $ws.Range($addr).NumberFormat = "#\ ##0"
    $ws.Range($addr).Font.Name = "Arial"
    $ws.Range($addr).Font.Size = 10
    $ws.Range($addr).Font.Color = -16777216
    $ws.Range($addr).Interior.ThemeColor = 0
}

# ---- Row 5: "disabilities Persons" + real data (was merged note row) ----
$ws.Range("A5:H5").UnMerge()
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.Color = -16777216
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").Font.Underline = 0
$ws.Range("A5").Interior.ThemeColor = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 21

$data5 = @(260, 228, 230, 236, 234, 238, 245, 245)
for ($i = 0; $i -lt 8; $i++) {
    $addr = "$($cols[$i])5"
    $ws.Range($addr).Value = $data5[$i]
    $ws.Range($addr).NumberFormat = "#\ ##0"
    $ws.Range($addr).Font.Name = "Arial"
    $ws.Range($addr).Font.Size = 10
    $ws.Range($addr).Font.Color = -16777216
    $ws.Range($addr).Interior.ThemeColor = 0
}
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2

# ---- Row 6: Source attribution (merged A6:H6) ----
$ws.Range("B6:H6").ClearContents()
$ws.Range("A6").Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$ws.Range("A6").Characters(1, 7).Font.Bold = $true
$ws.Range("A6").Characters(1, 7).Font.Underline = 2
$ws.Range("A6:H6").Font.Name = "Arial"
$ws.Range("A6:H6").Font.Size = 9
$ws.Range("A6:H6").Font.Color = -16777216
$ws.Range("A6:H6").Interior.ThemeColor = 0
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true
$ws.Range("A6:H6").Merge()
$ws.Rows.Item(6).RowHeight = 27.75

# ---- Selection matching the saved view state ----
$ws.Range("A1:I1").Select()
